$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.005.82"
$ws.Range("E2").Value = "  -0.21%  "
$ws.Range("D3").Value = "1.870.80"
$ws.Range("E3").Value = "  -0.11%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'312.29"
$ws.Range("E5").Value = "  -0.40%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("D7").Value = "'0.5137"
$ws.Range("E7").Value = "  +1.87%  "
$ws.Range("D8").Value = "'0.3851"
$ws.Range("E8").Value = "  +0.35%  "
$ws.Range("D9").Value = "'0.08301"
$ws.Range("E9").Value = "  -3.48%  "
$ws.Range("D10").Value = "'1.112"
$ws.Range("E10").Value = "  -0.23%  "
$ws.Range("D11").Value = "'41.51"
$ws.Range("E11").Value = "  -0.08%  "
$ws.Range("D12").Value = "'6.207"
$ws.Range("E12").Value = "  -1.47%  "
$ws.Range("D13").Value = "'20.56"
$ws.Range("E13").Value = "  -0.35%  "
$ws.Range("D14").Value = "1.848.04"
$ws.Range("E14").Value = "  -1.43%  "
$ws.Range("D15").Value = "'7.282"
$ws.Range("E15").Value = "  +1.32%  "
$ws.Range("D16").Value = "'1.002"
$ws.Range("E16").Value = "  -0.15%  "
$ws.Range("D17").Value = "'0.00001097"
$ws.Range("E17").Value = "  -0.03%  "
$ws.Range("D18").Value = "'90.66"
$ws.Range("E18").Value = "  -0.15%  "
$ws.Range("D19").Value = "'0.06657"
$ws.Range("E19").Value = "  +0.47%  "
$ws.Range("D20").Value = "'17.71"
$ws.Range("E20").Value = "  -1.74%  "
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("D22").Value = "'6.022"
$ws.Range("E22").Value = "  -0.84%  "
$ws.Range("D23").Value = "28.030.27"
$ws.Range("E23").Value = "  -0.26%  "
$ws.Range("D24").Value = "'11.09"
$ws.Range("E24").Value = "  -2.72%  "
$ws.Range("D25").Value = "'2.246"
$ws.Range("E25").Value = "  -0.86%  "
$ws.Range("D26").Value = "2.069.91"
$ws.Range("E26").Value = "  -1.04%  "
$ws.Range("D27").Value = "'2.516"
$ws.Range("E27").Value = "  -2.36%  "
$ws.Range("D28").Value = "'157.42"
$ws.Range("E28").Value = "  +0.39%  "
$ws.Range("D29").Value = "'20.59"
$ws.Range("E29").Value = "  -0.54%  "
$ws.Range("D30").Value = "'125.48"
$ws.Range("E30").Value = "  -0.58%  "
$ws.Range("D31").Value = "'0.1066"
$ws.Range("E31").Value = "  +1.26%  "
$ws.Range("D32").Value = "'1.032"
$ws.Range("E32").Value = "  -2.45%  "
$ws.Range("D33").Value = "'5.837"
$ws.Range("E33").Value = "  +4.24%  "
$ws.Range("D34").Value = "'3.597"
$ws.Range("E34").Value = "  -0.09%  "
$ws.Range("D35").Value = "'9.446"
$ws.Range("E35").Value = "  -1.52%  "
$ws.Range("D36").Value = "'0.02415"
$ws.Range("E36").Value = "  -0.79%  "
$ws.Range("D37").Value = "'0.06515"
$ws.Range("E37").Value = "  -0.86%  "
$ws.Range("D38").Value = "'0.2195"
$ws.Range("E38").Value = "  +1.11%  "
$ws.Range("D39").Value = "'0.6556"
$ws.Range("E39").Value = "  +2.95%  "
$ws.Range("D40").Value = "'1.202"
$ws.Range("E40").Value = "  -0.57%  "
$ws.Range("D41").Value = "'5.006"
$ws.Range("E41").Value = "  +2.54%  "
$ws.Range("D42").Value = "'1.209"
$ws.Range("E42").Value = "  -2.37%  "
$ws.Range("D43").Value = "'11.21"
$ws.Range("E43").Value = "  -2.51%  "
$ws.Range("D44").Value = "'0.6126"
$ws.Range("E44").Value = "  +2.28%  "
$ws.Range("D45").Value = "'13.00"
$ws.Range("E45").Value = "  -0.94%  "
$ws.Range("D46").Value = "'1.278"
$ws.Range("E46").Value = "  -0.39%  "
$ws.Range("D47").Value = "'3.673"
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("D48").Value = "'2.019"
$ws.Range("E48").Value = "  +1.65%  "
$ws.Range("D49").Value = "'1.214"
$ws.Range("E49").Value = "  -1.04%  "
$ws.Range("D50").Value = "'120.93"
$ws.Range("E50").Value = "  -0.41%  "
$ws.Range("D51").Value = "'78.02"
$ws.Range("E51").Value = "  -2.66%  "
